$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear RM 8's F value (F3) -> becomes blank
$ws.Range("F3").Value = $null

# Remove the "RM 232" row (row 26) entirely - all following rows shift up
$ws.Rows.Item(26).Delete()

# After the above deletion, "SC 92" (originally row 28) is now row 27.
# Remove it too - all following rows shift up again.
$ws.Rows.Item(27).Delete()

# Fill in / clear some previously-missing values now that rows have shifted:
# "SC 5" is now row 26 -> set its C (B column header) value
$ws.Range("C26").Value = 10.8

# "SC 101" is now row 27 -> clear its C value (now missing)
$ws.Range("C27").Value = $null

# "SC 232" is now the last row (33) -> fill in its C and F values
$ws.Range("C33").Value = 10.4
$ws.Range("F33").Value = 17.53
